# Revenues sheet: add a numeric "Release Month" column (B) right after the
# month-name column, shifting the existing "Global Box Office ($)" columns
# one slot to the right, then re-sort the data rows into calendar order
# (Jan..Dec) using the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column B. This shifts the old B ("Global Box Office ($)")
# to C and the old C ("Global Box Office ($)") to D, and carries over column
# A's style (s="1") to the new column - matching the target formatting.
$ws.Columns("B").Insert()

# Header for the newly inserted column.
$ws.Cells.Item(1, 2).Value = "Release Month"

# Populate the month-number (1-12) for every data row based on the
# three-letter month name already present in column A.
$monthNumbers = @{
    "Jan" = 1;  "Feb" = 2;  "Mar" = 3;  "Apr" = 4;
    "May" = 5;  "Jun" = 6;  "Jul" = 7;  "Aug" = 8;
    "Sep" = 9;  "Oct" = 10; "Nov" = 11; "Dec" = 12
}

$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $monthName = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Value = $monthNumbers[$monthName]
}

# Re-sort the data (rows 2..last) into ascending calendar order using the
# new "Release Month" numeric column as the sort key.
$dataRange = $ws.Range("A1:D$lastRow")
$sortKey = $ws.Range("B2:B$lastRow")
$dataRange.Sort($sortKey, 1, $null, $null, 1, $null, 1, 1)
